# Applies the two changes captured in the target OOXML diff:
#   1. The table on slide 16 gets a new table style (tableStyleId swap).
#   2. The presentation's (slide-master) theme color scheme is swapped
#      from the custom "Integral" palette to the stock "Office Theme"
#      palette (theme1.xml <-> theme2.xml content swap as seen upstream).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 ------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{13CC0A90-3358-466B-94BD-3D0CB17A9108}")
    }
}

# --- 2. Swap the theme colors to the "Office Theme" palette ---------------
function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Order matches ColorScheme.Colors(1..12):
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$scheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $scheme.Colors($i).RGB = HexToComRgb($officeThemeColors[$i - 1])
}
